# "Chansey" sheet update: add a new "graduation" dialogue block (2 new rows)
# after the existing table, and promote row 11 (the last row of the
# previous block) to a "block boundary" row by giving it a thin top+bottom
# border — matching the same visual convention already used by rows 6 and 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to build literal backslash-escape sequences (\n, \') that must be
# written verbatim into the cell text (this workbook stores raw "\n"/"\'"
# two-character sequences in its strings, not real newlines/quotes).
$BS = [char]92
$AP = [char]39
function Lit([string]$s) {
    return $s.Replace('~AP~', "$BS$AP")
}

# --- Row 11: turn it into a boundary row (thin top + thin bottom border) ---
# Clone the formatting already used for boundary rows (row 6) onto row 11,
# then thin the medium border down to a thin one so it matches the new
# border style used only for this row.
$ws.Range("A6:E6").Copy()
$ws.Range("A11:E11").PasteSpecial(-4122)   # xlPasteFormats
$rng11 = $ws.Range("A11:E11")
$rng11.Borders.Item(8).Weight = 2          # xlEdgeTop    -> xlThin
$rng11.Borders.Item(9).Weight = 2          # xlEdgeBottom -> xlThin

# --- Rows 12 & 13: new dialogue entries ---
# Copy the plain (no border) row formatting used by normal rows (row 9,
# columns B:E only -- column A is intentionally left untouched/empty here,
# same as rows 7-9) onto the two new rows.
$ws.Range("B9:E9").Copy()
$ws.Range("B12:E12").PasteSpecial(-4122)
$ws.Range("B13:E13").PasteSpecial(-4122)
$ws.Rows.Item(12).RowHeight = 21.6
$ws.Rows.Item(13).RowHeight = 21.6

$ws.Range("B12").Value2 = 137
$ws.Range("B13").Value2 = 140

# Column C = English, column D = Russian translation, column E = converted
# (mojibake-style re-encoded) string -- same pattern as every other row.
# Values are written in C,C,D,D,E,E order so the shared-string table gets
# populated in the same order the reference workbook uses.
$ws.Range("C12").Value2 = Lit " Eep! I~AP~m happy you graduated! ♪$($BS)nCongratulations!"
$ws.Range("C13").Value2 = Lit " I hope you~AP~ll keep up with your$($BS)nexploring! ♪"

$ws.Range("D12").Value2 = Lit " Иии! Я так рада, что вы$($BS)nвыпустились! ♪ Поздравляю!"
$ws.Range("D13").Value2 = Lit " Надеюсь, вы и дальше будете$($BS)nзаниматься исследованиями! ♪"

$ws.Range("E12").Value2 = Lit " Ééé! Ÿ óàë ñàäà, œóï âú$($BS)nâúðôòóéìéòû! ♪ Ðïèäñàâìÿý!"
$ws.Range("E13").Value2 = Lit " Îàäåýòû, âú é äàìûšå áôäåóå$($BS)nèàîéíàóûòÿ éòòìåäïâàîéÿíé! ♪"

# Match the recorded selection in the updated workbook.
$ws.Range("C6").Select() | Out-Null
